$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (date label + 12 numeric columns B:M)
$newData = @(
    @("27-09-2021", 4.19, 4.43, 4.6,  4.76, 5.22, -1.02, 0.29, 0.65, 1.01, 1.3,  1.84, 2.01),
    @("28-09-2021", 4.25, 4.5,  4.67, 4.83, 5.29, -1.03, 0.33, 0.6899999999999999, 1.06, 1.35, 1.9,  2.08),
    @("29-09-2021", 4.35, 4.6,  4.78, 4.95, 5.43, -1.16, 0.33, 0.73, 1.1,  1.41, 1.98, 2.15),
    @("30-09-2021", 4.38, 4.64, 4.84, 5.02, 5.47, -1.14, 0.36, 0.82, 1.18, 1.51, 2.09, 2.26),
    @("01-10-2021", 4.44, 4.71, 4.91, 5.09, 5.54, -0.92, 0.39, 0.85, 1.23, 1.57, 2.16, 2.34)
)

$startRow = 189
$endRow = $startRow + $newData.Count - 1

# Write the date labels through a formula first, so Excel's text parser
# doesn't reinterpret date-like strings (e.g. "01-10-2021") as a date
# serial number, then convert the formulas into plain static text values.
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $dateLabel = $newData[$i][0]
    $ws.Cells.Item($row, 1).Formula = '="' + $dateLabel + '"'
}

$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Fill in the numeric columns B:M
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $values = $newData[$i]

    for ($col = 2; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
